$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Replace the data table (rows 16-32) with the new dataset --------------
# Columns: B=Tipo Doc, C=N Doc Trabajador, D=Nombre Trabajador, E=Periodo Mora,
#          F=Valor Mora, G=Salario Basico

$data = @(
    @("CC","1047383187","ANTONIO JOSE BUELVAS CHAGUI","1709",36000,1350000),
    @("CC","1047383187","ANTONIO JOSE BUELVAS CHAGUI","1710",54000,1350000),
    @("CC","1047383187","ANTONIO JOSE BUELVAS CHAGUI","1711",54000,1350000),
    @("CC","1047383187","ANTONIO JOSE BUELVAS CHAGUI","1712",54000,1350000),
    @("CC","1047383187","ANTONIO JOSE BUELVAS CHAGUI","1801",54000,1350000),
    @("CC","1047383187","ANTONIO JOSE BUELVAS CHAGUI","1802",54000,1350000),
    @("CC","1047383187","ANTONIO JOSE BUELVAS CHAGUI","1803",54000,1350000),
    @("CC","1047383187","ANTONIO JOSE BUELVAS CHAGUI","1804",54000,1350000),
    @("CC","1047383187","ANTONIO JOSE BUELVAS CHAGUI","1805",54000,1350000),
    @("CC","1047383187","ANTONIO JOSE BUELVAS CHAGUI","1806",54000,1350000),
    @("CC","1047383187","ANTONIO JOSE BUELVAS CHAGUI","1807",54000,1350000),
    @("CC","1047383187","ANTONIO JOSE BUELVAS CHAGUI","1808",54000,1350000),
    @("CC","1047383187","ANTONIO JOSE BUELVAS CHAGUI","1809",54000,1350000),
    @("CC","1047383187","ANTONIO JOSE BUELVAS CHAGUI","1810",54000,1350000),
    @("CC","1047383187","ANTONIO JOSE BUELVAS CHAGUI","1811",54000,1350000),
    @("CC","1044930799","BALTAZAR DE JESUS POSADA LIGARDO","1903",7939,1190904),
    @("CC","45504211","DORYS VALENCIA BALLESTAS","2109",50384,6297977)
)

$row = 16
foreach ($r in $data) {
    $ws.Range("B$row").Value = $r[0]
    $ws.Range("C$row").Value = $r[1]
    $ws.Range("D$row").Value = $r[2]
    $ws.Range("E$row").Value = $r[3]
    $ws.Range("F$row").Value = $r[4]
    $ws.Range("G$row").Value = $r[5]
    $row++
}

# Row 32 (the new last data row) takes over the distinctive "closing" border
# formatting that used to live on row 40 before the table shrank.
$ws.Range("B40:J40").Copy()
$ws.Range("B32").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The old rows 33-44 (leftover data rows + blank gap) are no longer needed;
# deleting them shifts the signature-block footer rows up from 45/46 to 37/38.
$ws.Rows("33:40").Delete()

# --- Update the summary header cells ---------------------------------------
$ws.Range("E11").Value = 850323     # VALOR MORA total
$ws.Range("C13").Value = 3          # Cant. Trabajadores
$ws.Range("F13").Value = 17         # Cant. Periodos
